$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 324586.4
$ws.Range("I80").Value = 613.2857
$ws.Range("J80").Value = 778148.8
$ws.Range("K80").Value = 1839.8571
$ws.Range("L80").Value = 2334446.4
$ws.Range("M80").Value = -841.8571000000002
$ws.Range("N80").Value = -2336442.4
$ws.Range("H83").Value = 324586.4
$ws.Range("I83").Value = 613.2857
$ws.Range("J83").Value = 778148.8
$ws.Range("K83").Value = 5519.571300000001
$ws.Range("L83").Value = 7003339.2
$ws.Range("M83").Value = -527.5713000000005
$ws.Range("N83").Value = -7013323.2
$ws.Range("H115").Value = 346.25
$ws.Range("I115").Value = 335
$ws.Range("K115").Value = 1005
$ws.Range("M115").Value = 562
$ws.Range("H132").Value = 1922.8788
$ws.Range("I132").Value = 2021.8334
$ws.Range("J132").Value = 933.3333
$ws.Range("K132").Value = 6065.5002
$ws.Range("L132").Value = 2799.9999
$ws.Range("M132").Value = -3535.5002
$ws.Range("N132").Value = -7859.9999
$ws.Range("H135").Value = 961.1277
$ws.Range("I135").Value = 538.0789
$ws.Range("J135").Value = 2747.3333
$ws.Range("K135").Value = 4842.7101
$ws.Range("L135").Value = 24725.9997
$ws.Range("M135").Value = -2307.7101
$ws.Range("N135").Value = -29795.9997
$ws.Range("H137").Value = 789.2
$ws.Range("I137").Value = 735.6896400000001
$ws.Range("J137").Value = 930.2727
$ws.Range("K137").Value = 2207.06892
$ws.Range("L137").Value = 2790.8181
$ws.Range("M137").Value = 342.9310799999998
$ws.Range("N137").Value = -7890.8181
$ws.Range("H138").Value = 984.21
$ws.Range("I138").Value = 539.2143
$ws.Range("J138").Value = 2022.5333
$ws.Range("K138").Value = 1617.6429
$ws.Range("L138").Value = 6067.5999
$ws.Range("M138").Value = 3522.3571
$ws.Range("N138").Value = -16347.5999
$ws.Range("H141").Value = 1991.3334
$ws.Range("I141").Value = 646.975
$ws.Range("J141").Value = 8713.125
$ws.Range("K141").Value = 1940.925
$ws.Range("L141").Value = 26139.375
$ws.Range("M141").Value = 3239.075
$ws.Range("N141").Value = -36499.375

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19615.365
$ws.Range("I32").Value = 19469.727
$ws.Range("K32").Value = 19469.727
$ws.Range("M32").Value = -19182.727
$ws.Range("H61").Value = 1076
$ws.Range("I61").Value = 708.17645
$ws.Range("J61").Value = 2639.25
$ws.Range("K61").Value = 708.17645
$ws.Range("L61").Value = 2639.25
$ws.Range("M61").Value = -496.17645
$ws.Range("N61").Value = -3063.25
$ws.Range("H74").Value = 747.92206
$ws.Range("I74").Value = 678.36664
$ws.Range("K74").Value = 678.36664
$ws.Range("M74").Value = 195.63336
$ws.Range("H77").Value = 747.92206
$ws.Range("I77").Value = 678.36664
$ws.Range("K77").Value = 3391.8332
$ws.Range("M77").Value = 976.1668
$ws.Range("H132").Value = 1267.2142
$ws.Range("I132").Value = 906.725
$ws.Range("J132").Value = 2168.4375
$ws.Range("K132").Value = 2720.175
$ws.Range("L132").Value = 6505.3125
$ws.Range("M132").Value = -190.1750000000002
$ws.Range("N132").Value = -11565.3125
$ws.Range("H136").Value = 1076
$ws.Range("I136").Value = 708.17645
$ws.Range("J136").Value = 2639.25
$ws.Range("K136").Value = 2124.52935
$ws.Range("L136").Value = 7917.75
$ws.Range("M136").Value = 425.4706499999998
$ws.Range("N136").Value = -13017.75

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 30685.584
$ws.Range("I20").Value = 101624.5
$ws.Range("J20").Value = 7039.278
$ws.Range("K20").Value = 101624.5
$ws.Range("L20").Value = 7039.278
$ws.Range("M20").Value = -101377.5
$ws.Range("N20").Value = -7533.278
$ws.Range("H64").Value = 1429.1111
$ws.Range("I64").Value = 1767.3334
$ws.Range("J64").Value = 1090.8889
$ws.Range("K64").Value = 1767.3334
$ws.Range("L64").Value = 1090.8889
$ws.Range("M64").Value = -1542.3334
$ws.Range("N64").Value = -1540.8889
$ws.Range("H67").Value = 1429.1111
$ws.Range("I67").Value = 1767.3334
$ws.Range("J67").Value = 1090.8889
$ws.Range("K67").Value = 1767.3334
$ws.Range("L67").Value = 1090.8889
$ws.Range("M67").Value = -987.3334
$ws.Range("N67").Value = -2650.8889
$ws.Range("H134").Value = 14637.25
$ws.Range("I134").Value = 1129.3969
$ws.Range("K134").Value = 3388.1907
$ws.Range("M134").Value = -853.1907000000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2283.9482
$ws.Range("I31").Value = 2232.2856
$ws.Range("J31").Value = 2419.5625
$ws.Range("K31").Value = 2232.2856
$ws.Range("L31").Value = 2419.5625
$ws.Range("M31").Value = -1937.2856
$ws.Range("N31").Value = -3009.5625
$ws.Range("H34").Value = 2283.9482
$ws.Range("I34").Value = 2232.2856
$ws.Range("J34").Value = 2419.5625
$ws.Range("K34").Value = 2232.2856
$ws.Range("L34").Value = 2419.5625
$ws.Range("M34").Value = -2030.2856
$ws.Range("N34").Value = -2823.5625
$ws.Range("H58").Value = 3016.283
$ws.Range("I58").Value = 871.25
$ws.Range("J58").Value = 13503.111
$ws.Range("K58").Value = 871.25
$ws.Range("L58").Value = 13503.111
$ws.Range("M58").Value = -668.25
$ws.Range("N58").Value = -13909.111
$ws.Range("H92").Value = 21666.666
$ws.Range("J92").Value = 21666.666
$ws.Range("L92").Value = 21666.666
$ws.Range("N92").Value = -26658.666
$ws.Range("H132").Value = 1020.5227
$ws.Range("I132").Value = 732.0909
$ws.Range("J132").Value = 1885.8182
$ws.Range("K132").Value = 2196.2727
$ws.Range("L132").Value = 5657.4546
$ws.Range("M132").Value = 333.7273
$ws.Range("N132").Value = -10717.4546
$ws.Range("H136").Value = 3016.283
$ws.Range("I136").Value = 871.25
$ws.Range("J136").Value = 13503.111
$ws.Range("K136").Value = 2613.75
$ws.Range("L136").Value = 40509.333
$ws.Range("M136").Value = -63.75
$ws.Range("N136").Value = -45609.333

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 5000
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 5000
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 15000
$ws.Range("M76").Value = $null
$ws.Range("N76").Value = -15766
$ws.Range("H79").Value = 5000
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 5000
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 15000
$ws.Range("M79").Value = $null
$ws.Range("N79").Value = -17652
$ws.Range("H113").Value = 709.5454999999999
$ws.Range("I113").Value = 1240
$ws.Range("J113").Value = 553.5294
$ws.Range("K113").Value = 3720
$ws.Range("L113").Value = 1660.5882
$ws.Range("M113").Value = -1550
$ws.Range("N113").Value = -6000.5882
$ws.Range("H122").Value = 1000970.3
$ws.Range("J122").Value = 1112100.4
$ws.Range("L122").Value = 10008903.6
$ws.Range("N122").Value = -10013803.6
$ws.Range("H131").Value = 10682837
$ws.Range("I131").Value = 83502070
$ws.Range("J131").Value = 26363.61
$ws.Range("K131").Value = 250506210
$ws.Range("L131").Value = 79090.83
$ws.Range("M131").Value = -250501170
$ws.Range("N131").Value = -89170.83

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3133.1538
$ws.Range("I126").Value = 3214.625
$ws.Range("J126").Value = 3002.8
$ws.Range("K126").Value = 9643.875
$ws.Range("L126").Value = 9008.400000000001
$ws.Range("M126").Value = -7173.875
$ws.Range("N126").Value = -13948.4
$ws.Range("H132").Value = 2191.08
$ws.Range("I132").Value = 2120.8948
$ws.Range("J132").Value = 2413.3333
$ws.Range("K132").Value = 6362.6844
$ws.Range("L132").Value = 7239.999899999999
$ws.Range("M132").Value = -3832.6844
$ws.Range("N132").Value = -12299.9999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4100.2856
$ws.Range("I7").Value = 4414.857
$ws.Range("J7").Value = 3785.7144
$ws.Range("K7").Value = 4414.857
$ws.Range("L7").Value = 3785.7144
$ws.Range("M7").Value = -4302.857
$ws.Range("N7").Value = -4009.7144
$ws.Range("H46").Value = 1412.7273
$ws.Range("I46").Value = 1713
$ws.Range("J46").Value = 887.25
$ws.Range("K46").Value = 1713
$ws.Range("L46").Value = 887.25
$ws.Range("M46").Value = -1525
$ws.Range("N46").Value = -1263.25
$ws.Range("H126").Value = 4100.2856
$ws.Range("I126").Value = 4414.857
$ws.Range("J126").Value = 3785.7144
$ws.Range("K126").Value = 13244.571
$ws.Range("L126").Value = 11357.1432
$ws.Range("M126").Value = -10774.571
$ws.Range("N126").Value = -16297.1432
$ws.Range("H132").Value = 1132.0641
$ws.Range("I132").Value = 1021.9315
$ws.Range("K132").Value = 3065.7945
$ws.Range("M132").Value = -535.7945
$ws.Range("H136").Value = 1715.6863
$ws.Range("I136").Value = 935.0465
$ws.Range("J136").Value = 5911.625
$ws.Range("K136").Value = 2805.1395
$ws.Range("L136").Value = 17734.875
$ws.Range("M136").Value = -255.1395000000002
$ws.Range("N136").Value = -22834.875

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 838.4426
$ws.Range("I132").Value = 562.21155
$ws.Range("J132").Value = 2434.4443
$ws.Range("K132").Value = 1686.63465
$ws.Range("L132").Value = 7303.3329
$ws.Range("M132").Value = 843.36535
$ws.Range("N132").Value = -12363.3329
$ws.Range("H136").Value = 286.3
$ws.Range("I136").Value = 186.67273
$ws.Range("J136").Value = 1382.2
$ws.Range("K136").Value = 560.01819
$ws.Range("L136").Value = 4146.6
$ws.Range("M136").Value = 1989.98181
$ws.Range("N136").Value = -9246.6
